$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 3409.1
$ws.Range("I98").Value = 913.2857
$ws.Range("J98").Value = 9232.666999999999
$ws.Range("K98").Value = 913.2857
$ws.Range("L98").Value = 9232.666999999999
$ws.Range("M98").Value = 584.7143
$ws.Range("N98").Value = -12228.667
# Row 111
$ws.Range("H111").Value = 594.0833
$ws.Range("I111").Value = 518.8570999999999
$ws.Range("K111").Value = 1556.5713
$ws.Range("M111").Value = 1510.4287
# Row 112
$ws.Range("H112").Value = 2988.9038
$ws.Range("J112").Value = 3092.3877
$ws.Range("L112").Value = 9277.163100000002
$ws.Range("N112").Value = -11493.1631
# Row 122
$ws.Range("H122").Value = 3409.1
$ws.Range("I122").Value = 913.2857
$ws.Range("J122").Value = 9232.666999999999
$ws.Range("K122").Value = 2739.8571
$ws.Range("L122").Value = 27698.001
$ws.Range("M122").Value = -289.8571000000002
$ws.Range("N122").Value = -32598.001
# Row 125
$ws.Range("H125").Value = 17330.334
$ws.Range("I125").Value = 20997
$ws.Range("K125").Value = 188973
$ws.Range("M125").Value = -186513
# Row 138
$ws.Range("H138").Value = 2670.5571
$ws.Range("J138").Value = 3118.302
$ws.Range("L138").Value = 9354.906000000001
$ws.Range("N138").Value = -19634.906

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5080.7637
$ws.Range("I32").Value = 1001.06976
$ws.Range("K32").Value = 1001.06976
$ws.Range("M32").Value = -714.06976
# Row 61
$ws.Range("H61").Value = 5117.4688
$ws.Range("I61").Value = 4176.407
$ws.Range("K61").Value = 4176.407
$ws.Range("M61").Value = -3964.407
# Row 74
$ws.Range("H74").Value = 4070.625
$ws.Range("I74").Value = 3727.2778
$ws.Range("J74").Value = 5100.6665
$ws.Range("K74").Value = 3727.2778
$ws.Range("L74").Value = 5100.6665
$ws.Range("M74").Value = -2853.2778
$ws.Range("N74").Value = -6848.6665
# Row 77
$ws.Range("H77").Value = 4070.625
$ws.Range("I77").Value = 3727.2778
$ws.Range("J77").Value = 5100.6665
$ws.Range("K77").Value = 18636.389
$ws.Range("L77").Value = 25503.3325
$ws.Range("M77").Value = -14268.389
$ws.Range("N77").Value = -34239.3325
# Row 110
$ws.Range("H110").Value = 2496.3333
$ws.Range("I110").Value = 2246.625
$ws.Range("K110").Value = 2246.625
$ws.Range("M110").Value = -201.625
# Row 122
$ws.Range("H122").Value = 4421.5107
$ws.Range("I122").Value = 3885.475
$ws.Range("K122").Value = 11656.425
$ws.Range("M122").Value = -9206.424999999999
# Row 132
$ws.Range("H132").Value = 2156.84
$ws.Range("I132").Value = 1955.0416
$ws.Range("K132").Value = 5865.1248
$ws.Range("M132").Value = -3335.1248
# Row 136
$ws.Range("H136").Value = 5117.4688
$ws.Range("I136").Value = 4176.407
$ws.Range("K136").Value = 12529.221
$ws.Range("M136").Value = -9979.221000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 58
$ws.Range("H58").Value = 27303.5
$ws.Range("J58").Value = 26898
$ws.Range("L58").Value = 26898
$ws.Range("N58").Value = -27486
# Row 105
$ws.Range("H105").Value = 3140.7273
$ws.Range("I105").Value = 2899.3
$ws.Range("K105").Value = 2899.3
$ws.Range("M105").Value = -1152.3
# Row 107
$ws.Range("H107").Value = 2629.8
$ws.Range("I107").Value = 2678.682
$ws.Range("K107").Value = 2678.682
$ws.Range("M107").Value = -758.6819999999998
# Row 134
$ws.Range("H134").Value = 2936.04
$ws.Range("I134").Value = 2936.04
$ws.Range("K134").Value = 8808.119999999999
$ws.Range("M134").Value = -6273.119999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5949.9214
$ws.Range("I31").Value = 6194.45
$ws.Range("K31").Value = 6194.45
$ws.Range("M31").Value = -5899.45
# Row 34
$ws.Range("H34").Value = 5949.9214
$ws.Range("I34").Value = 6194.45
$ws.Range("K34").Value = 6194.45
$ws.Range("M34").Value = -5992.45
# Row 99
$ws.Range("H99").Value = 4998.25
$ws.Range("I99").Value = 4498.3
$ws.Range("K99").Value = 4498.3
$ws.Range("M99").Value = -3000.3
# Row 122
$ws.Range("H122").Value = 94283.03
$ws.Range("I122").Value = 123409.52
$ws.Range("J122").Value = 3262.75
$ws.Range("K122").Value = 370228.56
$ws.Range("L122").Value = 9788.25
$ws.Range("M122").Value = -367778.56
$ws.Range("N122").Value = -14688.25
# Row 126
$ws.Range("H126").Value = 4998.25
$ws.Range("I126").Value = 4498.3
$ws.Range("K126").Value = 13494.9
$ws.Range("M126").Value = -11024.9
# Row 132
$ws.Range("H132").Value = 1516.7826
$ws.Range("I132").Value = 1176.9412
$ws.Range("K132").Value = 3530.8236
$ws.Range("M132").Value = -1000.8236
# Row 134
$ws.Range("H134").Value = 906.95654
$ws.Range("I134").Value = 819.15
$ws.Range("K134").Value = 2457.45
$ws.Range("M134").Value = 77.55000000000018

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 152564400
$ws.Range("J4").Value = 20350000
$ws.Range("L4").Value = 61050000
$ws.Range("N4").Value = -61050224
# Row 107
$ws.Range("H107").Value = 503.65
$ws.Range("I107").Value = 474.33334
$ws.Range("K107").Value = 1423.00002
$ws.Range("M107").Value = 496.9999800000001
# Row 136
$ws.Range("H136").Value = 6172.6
$ws.Range("I136").Value = 5236.0527
$ws.Range("K136").Value = 15708.1581
$ws.Range("M136").Value = -10608.1581
# Row 140
$ws.Range("H140").Value = 1253.1052
$ws.Range("I140").Value = 1181.0555
$ws.Range("K140").Value = 3543.1665
$ws.Range("M140").Value = 1636.8335

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 300321.5
$ws.Range("J95").Value = 300321.5
$ws.Range("L95").Value = 300321.5
$ws.Range("N95").Value = -305813.5
# Row 113
$ws.Range("H113").Value = 1911.125
$ws.Range("I113").Value = 1911.125
$ws.Range("K113").Value = 1911.125
$ws.Range("M113").Value = 258.875
# Row 122
$ws.Range("H122").Value = 7080
$ws.Range("J122").Value = 7075.1665
$ws.Range("L122").Value = 21225.4995
$ws.Range("N122").Value = -26125.4995
# Row 126
$ws.Range("H126").Value = 8020.773
$ws.Range("J126").Value = 10163.111
$ws.Range("L126").Value = 30489.333
$ws.Range("N126").Value = -35429.333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1928.1428
$ws.Range("I100").Value = 1932.8334
$ws.Range("K100").Value = 1932.8334
$ws.Range("M100").Value = -1391.8334
# Row 132
$ws.Range("H132").Value = 34966.668
$ws.Range("I132").Value = 51250
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 153750
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -151220
$ws.Range("N132").Value = -12260

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 324
$ws.Range("I107").Value = 324
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 972
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 948
$ws.Range("N107").ClearContents()
# Row 113
$ws.Range("H113").Value = 2055.1667
$ws.Range("I113").Value = 1024.7222
$ws.Range("J113").Value = 5146.5
$ws.Range("K113").Value = 3074.1666
$ws.Range("L113").Value = 15439.5
$ws.Range("M113").Value = -904.1665999999996
$ws.Range("N113").Value = -19779.5
# Row 132
$ws.Range("H132").Value = 5292.3726
$ws.Range("I132").Value = 3545.476
$ws.Range("J132").Value = 13444.556
$ws.Range("K132").Value = 10636.428
$ws.Range("L132").Value = 40333.66800000001
$ws.Range("M132").Value = -8106.428
$ws.Range("N132").Value = -45393.66800000001
